$wb = $excel.ActiveWorkbook

# New "shadow price to market price ratio" coefficient row, added to every
# cost sheet (infrastructure, mobility, time). Font color FF222222 = 2236962.
$rpcFontColor = 2236962

# --- infrastructure sheet: infrast_cost_rpc ------------------------------
$wsInfra = $wb.Worksheets.Item("infrastructure")
$wsInfra.Range("A15").Value = "infrast_cost_rpc"
$wsInfra.Range("B15").Style = "Normal"
$wsInfra.Range("B15").Value = 0.82
$wsInfra.Range("B15").Font.Color = $rpcFontColor
$wsInfra.Range("C15").Value = "Shadow price to market price ratio in infrastructure cost (coeff)."
$wsInfra.Activate()
$wsInfra.Range("A15:C15").Select()

# --- mobility sheet: mobility_cost_rpc -----------------------------------
$wsMob = $wb.Worksheets.Item("mobility")
$wsMob.Range("A26").Value = "mobility_cost_rpc"
$wsMob.Range("B26").Style = "Normal"
$wsMob.Range("B26").Value = 0.82
$wsMob.Range("B26").Font.Color = $rpcFontColor
$wsMob.Range("C26").Value = "Shadow price to market price ratio in mobility cost (coeff)."
$wsMob.Activate()
$wsMob.Range("B27").Select()

# --- time sheet: time_cost_rpc -------------------------------------------
$wsTime = $wb.Worksheets.Item("time")
$wsTime.Range("A7").Value = "time_cost_rpc"
$wsTime.Range("B7").Style = "Normal"
$wsTime.Range("B7").Value = 1
$wsTime.Range("B7").Font.Color = $rpcFontColor
$wsTime.Range("C7").Value = "Shadow price to market price ratio in time cost (coeff)."
$wsTime.Activate()
$wsTime.Range("C7").Select()
